$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues
$xlPasteValues = -4163

$ws.Range("D2").Value = "30.601.86"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "1.673.40"
$ws.Range("E3").Value = "  +2.53%  "
$ws.Range("D4").Formula = "=""0.998"""
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial($xlPasteValues)
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Formula = "=""219.85"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial($xlPasteValues)
$ws.Range("E5").Value = "  +2.59%  "
$ws.Range("D6").Formula = "=""0.528"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial($xlPasteValues)
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("D7").Formula = "=""0.998"""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial($xlPasteValues)
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +4.21%  "
$ws.Range("E9").Value = "  +2.83%  "
$ws.Range("D11").Formula = "=""0.0906"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial($xlPasteValues)
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "1.913.84"
$ws.Range("E12").Value = "  +2.53%  "
$ws.Range("E13").Value = "  +9.28%  "
$ws.Range("D14").Formula = "=""10.25"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial($xlPasteValues)
$ws.Range("E14").Value = "  +12.31%  "
$ws.Range("D15").Value = "1.664.74"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("E16").Value = "  +4.30%  "
$ws.Range("D17").Value = "30.599.88"
$ws.Range("D18").Formula = "=""66.46"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial($xlPasteValues)
$ws.Range("E18").Value = "  +3.86%  "
$ws.Range("D19").Formula = "=""243.27"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial($xlPasteValues)
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("E20").Value = "  +3.45%  "
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("E22").Value = "  +3.76%  "
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Formula = "=""158.64"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial($xlPasteValues)
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("E26").Value = "  +2.39%  "
$ws.Range("E27").Value = "  +2.74%  "
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("D29").Formula = "=""0.999"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial($xlPasteValues)
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("E31").Value = "  +3.17%  "
$ws.Range("E32").Value = "  +3.00%  "
$ws.Range("D33").Formula = "=""3.30"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial($xlPasteValues)
$ws.Range("E33").Value = "  +3.94%  "
$ws.Range("D34").Value = "1.488.20"
$ws.Range("E34").Value = "  +4.48%  "
$ws.Range("E35").Value = "  +7.04%  "
$ws.Range("D36").Formula = "=""84.94"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial($xlPasteValues)
$ws.Range("E36").Value = "  +12.52%  "
$ws.Range("D37").Formula = "=""1.02"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial($xlPasteValues)
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("E38").Value = "  +8.60%  "
$ws.Range("E39").Value = "  +5.38%  "
$ws.Range("E40").Value = "  -4.35%  "
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("D43").Formula = "=""0.0499"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial($xlPasteValues)
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Formula = "=""0.999"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial($xlPasteValues)
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Formula = "=""51.42"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial($xlPasteValues)
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("D48").Formula = "=""5.52"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial($xlPasteValues)
$ws.Range("E48").Value = "  +3.06%  "
$ws.Range("D49").Value = "1.806.32"
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("D50").Formula = "=""94.80"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial($xlPasteValues)
$ws.Range("E50").Value = "  +4.72%  "
$ws.Range("E51").Value = "  -0.55%  "

$excel.CutCopyMode = 0

